# Apply scheduled Tiamat_Profits profitability recalculation to Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Columns H..N hold currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) figures refreshed by the runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: "Glazed and Confused" (Clear Glass Lens)
$ws.Range("H33").Value = 637.6667
$ws.Range("I33").Value = 207.22223
$ws.Range("J33").Value = 1283.3334
$ws.Range("K33").Value = 207.22223
$ws.Range("L33").Value = 1283.3334
$ws.Range("M33").Value = 21.77777
$ws.Range("N33").Value = -1741.3334

# Row 64: "Forged from the Void" (Void Glue)
$ws.Range("H64").Value = 500000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 500000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 500000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -500496

# Row 67: "Dodging the Draft (L)" (Void Glue)
$ws.Range("H67").Value = 500000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 500000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 500000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -501716

# Row 76: "Warding Off Temptation" (Enchanted Hardsilver Ink)
$ws.Range("H76").Value = 33336528
$ws.Range("I76").Value = 33336528
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 33336528
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -33336213
$ws.Range("N76").ClearContents()

# Row 79: "The Garden of Arcane Delights (L)" (Enchanted Hardsilver Ink)
$ws.Range("H79").Value = 33336528
$ws.Range("I79").Value = 33336528
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 33336528
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -33335436
$ws.Range("N79").ClearContents()

# Row 129: "Practical Command" (Commanding Craftsman's Draught)
$ws.Range("H129").Value = 823.4545000000001
$ws.Range("I129").Value = 576.2
$ws.Range("J129").Value = 1029.5
$ws.Range("K129").Value = 1728.6
$ws.Range("L129").Value = 3088.5
$ws.Range("M129").Value = 3271.4
$ws.Range("N129").Value = -13088.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32: "Ingot We Trust" (Steel Ingot)
$ws.Range("H32").Value = 5309.56
$ws.Range("I32").Value = 4876.796
$ws.Range("J32").Value = 11059.143
$ws.Range("K32").Value = 4876.796
$ws.Range("L32").Value = 11059.143
$ws.Range("M32").Value = -4589.796
$ws.Range("N32").Value = -11633.143

# Row 34: "Insistent Sallets" (Steel Sallet)
$ws.Range("H34").Value = 65577.25
$ws.Range("I34").Value = 2225
$ws.Range("J34").Value = 86694.664
$ws.Range("K34").Value = 2225
$ws.Range("L34").Value = 86694.664
$ws.Range("M34").Value = -1954
$ws.Range("N34").Value = -87236.664

# Row 61: "Dealing with the Tough Stuff" (Cobalt Ingot)
$ws.Range("H61").Value = 4050.15
$ws.Range("I61").Value = 3147.2354
$ws.Range("J61").Value = 9166.666999999999
$ws.Range("K61").Value = 3147.2354
$ws.Range("L61").Value = 9166.666999999999
$ws.Range("M61").Value = -2935.2354
$ws.Range("N61").Value = -9590.666999999999

# Row 74: "As the Bolt Flies" (Titanium Nugget)
$ws.Range("H74").Value = 19997.777
$ws.Range("I74").Value = 28540.676
$ws.Range("J74").Value = 1404.4117
$ws.Range("K74").Value = 28540.676
$ws.Range("L74").Value = 1404.4117
$ws.Range("M74").Value = -27666.676
$ws.Range("N74").Value = -3152.4117

# Row 77: "Heavy Metal Banned (L)" (Titanium Nugget)
$ws.Range("H77").Value = 19997.777
$ws.Range("I77").Value = 28540.676
$ws.Range("J77").Value = 1404.4117
$ws.Range("K77").Value = 142703.38
$ws.Range("L77").Value = 7022.058500000001
$ws.Range("M77").Value = -138335.38
$ws.Range("N77").Value = -15758.0585

# Row 136: "Metal with Mettle" (Cobalt Tungsten Ingot)
$ws.Range("H136").Value = 4050.15
$ws.Range("I136").Value = 3147.2354
$ws.Range("J136").Value = 9166.666999999999
$ws.Range("K136").Value = 9441.706200000001
$ws.Range("L136").Value = 27500.001
$ws.Range("M136").Value = -6891.706200000001
$ws.Range("N136").Value = -32600.001

$ws = $wb.Worksheets.Item("BSM")
# Row 86: "Through Thick and Thin" (Adamantite Nugget)
$ws.Range("H86").Value = 242785.38
$ws.Range("I86").Value = 1482.2222
$ws.Range("J86").Value = 637645.0600000001
$ws.Range("K86").Value = 1482.2222
$ws.Range("L86").Value = 637645.0600000001
$ws.Range("M86").Value = -359.2221999999999
$ws.Range("N86").Value = -639891.0600000001

# Row 89: "Piercing Eyes Deserve Piercing Shafts (L)" (Adamantite Nugget)
$ws.Range("H89").Value = 242785.38
$ws.Range("I89").Value = 1482.2222
$ws.Range("J89").Value = 637645.0600000001
$ws.Range("K89").Value = 7411.111
$ws.Range("L89").Value = 3188225.3
$ws.Range("M89").Value = -1795.111
$ws.Range("N89").Value = -3199457.3

$ws = $wb.Worksheets.Item("CRP")
# Row 31: "Wall Not Found" (Walnut Lumber)
$ws.Range("H31").Value = 9874.654
$ws.Range("I31").Value = 6602.971
$ws.Range("K31").Value = 6602.971
$ws.Range("M31").Value = -6307.971

# Row 34: "Armoires of the Rich and Famous" (Walnut Lumber)
$ws.Range("H34").Value = 9874.654
$ws.Range("I34").Value = 6602.971
$ws.Range("K34").Value = 6602.971
$ws.Range("M34").Value = -6400.971

# Row 94: "Beech, Please" (Beech Lumber)
$ws.Range("H94").Value = 1948.0416
$ws.Range("I94").Value = 2446.2856
$ws.Range("J94").Value = 1742.8823
$ws.Range("K94").Value = 2446.2856
$ws.Range("L94").Value = 1742.8823
$ws.Range("M94").Value = -1995.2856
$ws.Range("N94").Value = -2644.8823

# Row 96: "Composition" (Larch Composite Bow)
$ws.Range("H96").Value = 21500
$ws.Range("J96").Value = 21500
$ws.Range("L96").Value = 21500
$ws.Range("N96").Value = -26992

$ws = $wb.Worksheets.Item("CUL")
# Row 81: "It Goes Down Smoothly" (Frozen Spirits)
$ws.Range("H81").Value = 1560.1428
$ws.Range("I81").Value = 670.3333
$ws.Range("J81").Value = 2227.5
$ws.Range("K81").Value = 2010.9999
$ws.Range("L81").Value = 6682.5
$ws.Range("M81").Value = -887.9999
$ws.Range("N81").Value = -8928.5

# Row 84: "Quenching the Flame (L)" (Frozen Spirits)
$ws.Range("H84").Value = 1560.1428
$ws.Range("I84").Value = 670.3333
$ws.Range("J84").Value = 2227.5
$ws.Range("K84").Value = 6032.9997
$ws.Range("L84").Value = 20047.5
$ws.Range("M84").Value = -416.9997000000003
$ws.Range("N84").Value = -31279.5

# Row 121: "A Cookie for Your Troubles" (Coffee Biscuit)
$ws.Range("H121").Value = 51742224
$ws.Range("I121").Value = 666
$ws.Range("J121").Value = 62521710
$ws.Range("K121").Value = 1998
$ws.Range("L121").Value = 187565130
$ws.Range("M121").Value = -688
$ws.Range("N121").Value = -187567750

# Row 131: "The Mountain Steeped" (Tsai tou Vounou)
$ws.Range("H131").Value = 19481284
$ws.Range("I131").Value = 336
$ws.Range("J131").Value = 23716272
$ws.Range("K131").Value = 1008
$ws.Range("L131").Value = 71148816
$ws.Range("M131").Value = 4032
$ws.Range("N131").Value = -71158896

$ws = $wb.Worksheets.Item("GSM")
# Row 33: "Thaumaturge Is Magic" (Fluorite Ring)
$ws.Range("H33").Value = 6019
$ws.Range("J33").Value = 6019
$ws.Range("L33").Value = 6019
$ws.Range("N33").Value = -6523

# Row 39: "One Man's Trash" (Horn Ring)
$ws.Range("H39").Value = 30000
$ws.Range("J39").Value = 30000
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -31064

# Row 97: "If I'd a Koppranickel for Every Time..." (Koppranickel Ingot)
$ws.Range("H97").Value = 866.86206
$ws.Range("I97").Value = 875.95
$ws.Range("J97").Value = 846.6667
$ws.Range("K97").Value = 875.95
$ws.Range("L97").Value = 846.6667
$ws.Range("M97").Value = -379.95
$ws.Range("N97").Value = -1838.6667

$ws = $wb.Worksheets.Item("LTW")
# Row 93: "Hide to Go Seek" (Gagana Leather)
$ws.Range("H93").Value = 1753.56
$ws.Range("I93").Value = 1205.2667
$ws.Range("J93").Value = 2576
$ws.Range("K93").Value = 1205.2667
$ws.Range("L93").Value = 2576
$ws.Range("M93").Value = 42.7333000000001
$ws.Range("N93").Value = -5072

# Row 133: "The Perfect Accessory" (Loboskin Amulet of Fending)
$ws.Range("H133").Value = 45414.445
$ws.Range("J133").Value = 45414.445
$ws.Range("L133").Value = 45414.445
$ws.Range("N133").Value = -50474.445

# Row 136: "Respect for Br'aax" (Br'aax Leather)
$ws.Range("H136").Value = 436321.9
$ws.Range("I136").Value = 667745.6
$ws.Range("J136").Value = 2402.5
$ws.Range("K136").Value = 2003236.8
$ws.Range("L136").Value = 7207.5
$ws.Range("M136").Value = -2000686.8
$ws.Range("N136").Value = -12307.5

$ws = $wb.Worksheets.Item("WVR")
# Row 81: "Where the Dragonflies, the Net Catches" (Crawler Silk)
$ws.Range("H81").Value = 1983.48
$ws.Range("I81").Value = 2153.2666
$ws.Range("K81").Value = 4306.5332
$ws.Range("M81").Value = -3245.5332

# Row 84: "To Kill a Dragon on Nameday (L)" (Crawler Silk)
$ws.Range("H84").Value = 1983.48
$ws.Range("I84").Value = 2153.2666
$ws.Range("K84").Value = 21532.666
$ws.Range("M84").Value = -16228.666

# Row 107: "Flax Wax" (Bright Linen Yarn)
$ws.Range("H107").Value = 346.07693
$ws.Range("I107").Value = 291.9
$ws.Range("J107").Value = 526.6667
$ws.Range("K107").Value = 875.6999999999999
$ws.Range("L107").Value = 1580.0001
$ws.Range("M107").Value = 1044.3
$ws.Range("N107").Value = -5420.0001
